# Apply "Model Updated to Function with FSN Design" changes to the muscle
# parameter table.
#
# For every data row (columns I through P), the meaning of the columns was
# changed:
#   I (Am)   -> new Am, scaled from the old Am by a constant factor (0.8335)
#   J (S)    -> replaced by a single new constant value for every row
#   K (Yoff) -> replaced by a single new constant value for every row
#   L (Lmin) -> replaced by a single new constant value for every row (-60 mV)
#   M (l0)   -> replaced by a single new constant value for every row (-40 mV)
#   N (Lmax) -> replaced by a single new constant value for every row (0 mV)
#   O (Lw)   -> now mirrors the new Am value (same as column I)
#   P (Fo)   -> now holds what used to be the Lmin value (old column L)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constants introduced by the new "Function with FSN Design" model.
$AmScale   = 0.8335
$NewS      = 921.04403669765156
$NewYoff   = -0.001
$NewLmin   = -60
$NewL0     = -40
$NewLmax   = 0

# Find the last populated data row (column A holds the muscle name).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    # Capture the "before" values we still need once the columns are
    # overwritten (old Am in column I, old Lmin in column L).
    $oldAm   = $ws.Cells.Item($r, 9).Value2    # column I
    $oldLmin = $ws.Cells.Item($r, 12).Value2   # column L

    $newAm = $oldAm * $AmScale

    $ws.Cells.Item($r, 9).Value2  = $newAm     # I: Am
    $ws.Cells.Item($r, 10).Value2 = $NewS       # J: S
    $ws.Cells.Item($r, 11).Value2 = $NewYoff    # K: Yoff
    $ws.Cells.Item($r, 12).Value2 = $NewLmin    # L: Lmin
    $ws.Cells.Item($r, 13).Value2 = $NewL0      # M: l0
    $ws.Cells.Item($r, 14).Value2 = $NewLmax    # N: Lmax
    $ws.Cells.Item($r, 15).Value2 = $newAm      # O: Lw (mirrors Am)
    $ws.Cells.Item($r, 16).Value2 = $oldLmin    # P: Fo (old Lmin)
}
